$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column retains its literal text representation
# (values like "1.000" or "0.000008033" would otherwise be
# auto-converted to numbers by Excel, losing trailing zeros / dot formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.778.98"
$ws.Range("E2").Value = "  -3.79%  "

$ws.Range("D3").Value = "1.820.00"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "279.32"

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("D7").Value = "0.5099"
$ws.Range("E7").Value = "  -4.20%  "

$ws.Range("D8").Value = "0.3546"
$ws.Range("E8").Value = "  -5.21%  "

$ws.Range("D9").Value = "44.50"
$ws.Range("E9").Value = "  -2.29%  "

$ws.Range("D10").Value = "0.06667"
$ws.Range("E10").Value = "  -7.16%  "

$ws.Range("D11").Value = "20.06"
$ws.Range("E11").Value = "  -7.18%  "

$ws.Range("D12").Value = "0.8279"
$ws.Range("E12").Value = "  -6.82%  "

$ws.Range("D13").Value = "0.07906"
$ws.Range("E13").Value = "  -3.13%  "

$ws.Range("D14").Value = "1.818.77"
$ws.Range("E14").Value = "  -4.14%  "

$ws.Range("D15").Value = "5.084"
$ws.Range("E15").Value = "  -4.04%  "

$ws.Range("D16").Value = "88.03"
$ws.Range("E16").Value = "  -5.23%  "

$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").Value = "14.12"
$ws.Range("E18").Value = "  -4.69%  "

$ws.Range("D19").Value = "0.000008033"
$ws.Range("E19").Value = "  -5.76%  "

$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.18%  "

$ws.Range("D21").Value = "25.823.63"
$ws.Range("E21").Value = "  -3.78%  "

$ws.Range("E22").Value = "  -4.50%  "

$ws.Range("D23").Value = "9.998"
$ws.Range("E23").Value = "  -5.84%  "

$ws.Range("D24").Value = "6.123"
$ws.Range("E24").Value = "  -4.15%  "

$ws.Range("D25").Value = "2.240"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").Value = "142.20"
$ws.Range("E26").Value = "  -2.98%  "

$ws.Range("D27").Value = "1.670"
$ws.Range("E27").Value = "  -3.89%  "

$ws.Range("D28").Value = "17.14"
$ws.Range("E28").Value = "  -4.94%  "

$ws.Range("D29").Value = "109.54"
$ws.Range("E29").Value = "  -3.91%  "

$ws.Range("D30").Value = "4.317"
$ws.Range("E30").Value = "  -8.34%  "

$ws.Range("D31").Value = "4.247"
$ws.Range("E31").Value = "  -7.82%  "

$ws.Range("D32").Value = "0.08753"
$ws.Range("E32").Value = "  -3.96%  "

$ws.Range("D33").Value = "0.04914"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("D34").Value = "0.7304"
$ws.Range("E34").Value = "  -10.00%  "

$ws.Range("D35").Value = "1.139"
$ws.Range("E35").Value = "  -2.47%  "

$ws.Range("D36").Value = "2.881"
$ws.Range("E36").Value = "  -2.79%  "

$ws.Range("D37").Value = "3.158"
$ws.Range("E37").Value = "  -1.39%  "

$ws.Range("D38").Value = "0.9993"
$ws.Range("E38").Value = "  -0.31%  "

$ws.Range("D39").Value = "2.363"
$ws.Range("E39").Value = "  -10.94%  "

$ws.Range("D40").Value = "0.01860"
$ws.Range("E40").Value = "  -4.86%  "

$ws.Range("E41").Value = "  -15.04%  "

$ws.Range("D42").Value = "0.9684"
$ws.Range("E42").Value = "  -9.21%  "

$ws.Range("D43").Value = "114.32"
$ws.Range("E43").Value = "  -0.49%  "

$ws.Range("D44").Value = "6.260"
$ws.Range("E44").Value = "  -4.30%  "

$ws.Range("D45").Value = "8.061"
$ws.Range("E45").Value = "  -8.78%  "

$ws.Range("D46").Value = "1.0000"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D47").Value = "0.4561"
$ws.Range("E47").Value = "  -11.73%  "

$ws.Range("D48").Value = "0.1371"
$ws.Range("E48").Value = "  -8.20%  "

$ws.Range("D49").Value = "36.62"
$ws.Range("E49").Value = "  -2.32%  "

$ws.Range("D50").Value = "9.210"
$ws.Range("E50").Value = "  -7.40%  "

$ws.Range("E51").Value = "  -8.39%  "
